$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.224.35"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "1.859.35"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'0.7154"
$ws.Range("D5").QuotePrefix = $false
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").Value = "'240.65"
$ws.Range("D6").QuotePrefix = $false
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.07757"
$ws.Range("D8").QuotePrefix = $false
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").Value = "'0.3077"
$ws.Range("D9").QuotePrefix = $false
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").Value = "'0.08258"
$ws.Range("D11").QuotePrefix = $false
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "1.859.78"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "'5.247"
$ws.Range("D13").QuotePrefix = $false
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "'90.33"
$ws.Range("D15").QuotePrefix = $false
$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").Value = "29.207.33"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "'5.870"
$ws.Range("D17").QuotePrefix = $false
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").Value = "'244.17"
$ws.Range("D18").QuotePrefix = $false
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("D19").Value = "'0.000007803"
$ws.Range("D19").QuotePrefix = $false
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").Value = "'13.17"
$ws.Range("D20").QuotePrefix = $false
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "2.104.04"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").QuotePrefix = $false
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'7.995"
$ws.Range("D23").QuotePrefix = $false
$ws.Range("E23").Value = "  +3.41%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").QuotePrefix = $false
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'0.1600"
$ws.Range("D25").QuotePrefix = $false
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("D26").Value = "'162.60"
$ws.Range("D26").QuotePrefix = $false
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "'8.926"
$ws.Range("D27").QuotePrefix = $false
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("D28").Value = "'18.28"
$ws.Range("D28").QuotePrefix = $false
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").Value = "'1.498"
$ws.Range("D29").QuotePrefix = $false
$ws.Range("E29").Value = "  +1.06%  "

$ws.Range("D30").Value = "'1.315"
$ws.Range("D30").QuotePrefix = $false
$ws.Range("E30").Value = "  -3.12%  "

$ws.Range("D31").Value = "'4.403"
$ws.Range("D31").QuotePrefix = $false
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("D32").Value = "'4.201"
$ws.Range("D32").QuotePrefix = $false
$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").Value = "'0.05191"
$ws.Range("D33").QuotePrefix = $false
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "'1.912"
$ws.Range("D34").QuotePrefix = $false
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").Value = "'1.172"
$ws.Range("D35").QuotePrefix = $false
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").Value = "'0.7275"
$ws.Range("D36").QuotePrefix = $false
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").Value = "'2.675"
$ws.Range("D37").QuotePrefix = $false
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "'0.01857"
$ws.Range("D38").QuotePrefix = $false
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("D39").Value = "'2.688"
$ws.Range("D39").QuotePrefix = $false
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").Value = "1.153.31"
$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("D41").Value = "'0.9042"
$ws.Range("D41").QuotePrefix = $false
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").Value = "'6.158"
$ws.Range("D42").QuotePrefix = $false
$ws.Range("E42").Value = "  +2.90%  "

$ws.Range("D43").Value = "'72.27"
$ws.Range("D43").QuotePrefix = $false
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").QuotePrefix = $false
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "'101.68"
$ws.Range("D45").QuotePrefix = $false
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "2.002.86"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("D47").Value = "'0.5219"
$ws.Range("D47").QuotePrefix = $false
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("D48").Value = "'1.769"
$ws.Range("D48").QuotePrefix = $false
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("D50").Value = "'9.321"
$ws.Range("D50").QuotePrefix = $false
$ws.Range("E50").Value = "  +1.88%  "

$ws.Range("D51").Value = "'2.870"
$ws.Range("D51").QuotePrefix = $false
$ws.Range("E51").Value = "  +1.46%  "

